$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "The initial problem ... river." -> "... river in one trip."
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "to the other side of the river.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "to the other side of the river in one trip.", 2)

# ---------------------------------------------------------------------------
# 2) Re-write the rest of that paragraph: "implies that it is equally
#    important ... leave any behind." becomes the new reasoning about
#    "in his absence" / two trips, while keeping the final sentence about
#    transporting everybody safely. (The bookmark is relocated separately
#    in step 3, right before "Ultimately".)
# ---------------------------------------------------------------------------
$oldPart = "implies that it is equally important to him to be able to grow those seeds" + `
    "; so it is possible that he is in a survival situation where food is scarce.  " + `
    "Ultimately the goal is to figure out how he can transport himself, the parrot, cat, and seeds" + `
    " safely and without having to leave any behind."

$newPart = "implies that they are equally important to him.  " + `
    "Because it states the words " + [char]0x201C + "in his absence" + [char]0x201D + ", I assume it means he " + `
    "has considered the possibility of having to take two trips already. " + `
    "Ultimately the goal is to figure out how he can transport himself, the parrot, cat, and seeds" + `
    " safely and without having to leave any behind.  "

$null = $d.Content.Find.Execute(
    $oldPart,
    $true, $false, $false, $false, $false, $true, 1, $false,
    $newPart, 2)

# ---------------------------------------------------------------------------
# 3) Move the "_GoBack" bookmark from the end of the "In breaking the
#    problem apart ..." paragraph to right before "Ultimately the goal..."
#    in the paragraph above.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
}

$locate = $d.Content
$null = $locate.Find.Execute("Ultimately the goal is to figure out", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bookmarkPoint = $d.Range($locate.Start, $locate.Start)
$d.Bookmarks.Add("_GoBack", $bookmarkPoint)

# ---------------------------------------------------------------------------
# 4) After the "In breaking the problem apart ... help him." paragraph,
#    add a blank paragraph followed by a new paragraph: "The sub-goals
#    then are to get off the riverbank"
# ---------------------------------------------------------------------------
$constraintsIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "In breaking the problem apart*help him.*") {
        $constraintsIndex = $i
        break
    }
}

$constraintsPara = $d.Paragraphs.Item($constraintsIndex)
$constraintsPara.Range.InsertParagraphAfter()

$blankPara = $d.Paragraphs.Item($constraintsIndex + 1)
$blankPara.Range.InsertParagraphAfter()

$subGoalsPara = $d.Paragraphs.Item($constraintsIndex + 2)
$subGoalsPara.Range.InsertAfter("The sub-goals then are to get off the riverbank")

Write-Output "done"
